$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
# "Volume 31   Number  21" -> "Volume 31   Number  22"
$ws.Range("A8").Value = "Volume 31   Number  22"
# "Report Covering the Week  5/20/2024  Through  5/26/2024"
#   -> "Report Covering the Week  5/27/2024  Through  6/2/2024"
$ws.Range("C9").Value = "Report Covering the Week  5/27/2024  Through  6/2/2024"

# --- Fix cell styles/types for cells whose underlying type changes ---
# Text ("---"/"***.*") -> Number: copy format from a same-row numeric donor cell
# first (so the destination keeps a plain numeric style), then the final value is
# written below in the generic value pass.
$ws.Range("F15").Copy($ws.Range("C15"))
$ws.Range("F27").Copy($ws.Range("C27"))
$ws.Range("F27").Copy($ws.Range("D27"))
$ws.Range("H27").Copy($ws.Range("E27"))

# Number -> Text: copy format+value directly from a same-row donor cell that
# already holds the literal text "0" - this sets both style and value in one shot,
# so F31 is intentionally NOT touched again in the value pass below.
$ws.Range("G31").Copy($ws.Range("F31"))

# --- Set the final values for every changed numeric cell ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 85.714285714285
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = 44.444444444444
$ws.Range("N15").Value = 18.181818181818
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 157
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = 34.188034188034
$ws.Range("L16").Value = 50.961538461538
$ws.Range("M16").Value = 9.790209790209
$ws.Range("N16").Value = -70.095238095238
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 22.5
$ws.Range("I17").Value = 233
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 16.5
$ws.Range("L17").Value = 28.021978021978
$ws.Range("M17").Value = 80.620155038759
$ws.Range("N17").Value = 42.944785276073
$ws.Range("C18").Value = 8
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 17.647058823529
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = 49.206349206349
$ws.Range("L18").Value = 46.875
$ws.Range("M18").Value = -31.884057971014
$ws.Range("N18").Value = -89.508928571428
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 362
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = 11.384615384615
$ws.Range("L19").Value = -13.397129186602
$ws.Range("M19").Value = 95.675675675675
$ws.Range("N19").Value = -34.892086330935
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 123
$ws.Range("J20").Value = 128
$ws.Range("K20").Value = -3.90625
$ws.Range("L20").Value = 21.782178217821
$ws.Range("M20").Value = 25.510204081632
$ws.Range("N20").Value = -86.498353457738
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 164
$ws.Range("H21").Value = 13.414634146341
$ws.Range("I21").Value = 982
$ws.Range("J21").Value = 840
$ws.Range("K21").Value = 16.904761904761
$ws.Range("L21").Value = 11.211778029445
$ws.Range("M21").Value = 39.687055476529
$ws.Range("N21").Value = -68.002606712284
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -35.714285714285
$ws.Range("I22").Value = 29
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = -39.583333333333
$ws.Range("L22").Value = 61.111111111111
$ws.Range("M22").Value = 93.333333333333
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 6.25
$ws.Range("F24").Value = 139
$ws.Range("G24").Value = 139
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 965
$ws.Range("J24").Value = 859
$ws.Range("K24").Value = 12.339930151338
$ws.Range("L24").Value = 28.324468085106
$ws.Range("M24").Value = 117.342342342342
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -21.739130434782
$ws.Range("G25").Value = 77
$ws.Range("H25").Value = -12.987012987013
$ws.Range("I25").Value = 541
$ws.Range("J25").Value = 463
$ws.Range("K25").Value = 16.846652267818
$ws.Range("L25").Value = 26.995305164319
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = 52.941176470588
$ws.Range("F26").Value = 113
$ws.Range("G26").Value = 73
$ws.Range("H26").Value = 54.794520547945
$ws.Range("I26").Value = 528
$ws.Range("J26").Value = 368
$ws.Range("K26").Value = 43.478260869565
$ws.Range("L26").Value = 37.142857142857
$ws.Range("M26").Value = 37.5
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 21.052631578947
$ws.Range("L27").Value = 15
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -40
$ws.Range("F28").Value = 11
$ws.Range("H28").Value = -35.294117647058
$ws.Range("I28").Value = 61
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = -1.612903225806
$ws.Range("L28").Value = 48.780487804878
# F31 already set above via Copy() from its text donor cell
